$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "Record"
$ws.Range("B18").Value = "Balanço Geral"
$ws.Range("C18").Value = "Transporte"
$ws.Range("D18").Value = "2025-04-01T12:12"
$ws.Range("E18").Value = "Negativo"
$ws.Range("F18").Value = "Permissionários de van do setor C fazem nova paralisação. Repórter *ao vivo* do Centro. Essa é a segunda paralisação do setor C, que faz a linha de localidades mais distantes. Passageiros pagam R$ 3,50 e o valor é complementado pela prefeitura. Permissionários na frente da prefeitura. Eles querem falar com representante da prefeitura e prefeito Wladimir Garotinho. Eles alegam que teve o repasse ontem à tarde, mas não receberam o valor devido de acordo com a bilhetagem. Valor estaria com inconsistências. Receberam de 80% a 90% a menos do que deveriam receber. Pediram através do judiciário que atende aos permissionários para que fossem atendidos e saber o que realmente está acontecendo. Sistema de bilhetagem foi implantado a pedido da prefeitura em 2019.Alegam que não tem como trabalhar desse jeito. Entrevista com passageiros. Entrevista com permissionário Jefferson Oliveira. *Com nota do IMTT*"
